# Freeze ejecutado: roster y freeze_time actualizados
# Apply roster updates: players swapped out, remaining rows shifted up,
# and newly added players placed into the vacated roster slots.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RosterRow {
    param(
        [int]$Row,
        [double]$PlayerId,
        [string]$PlayerName,
        [string]$ProTeam,
        [string]$LineupSlot,
        [string]$NbaPlayerId
    )

    $ws.Range("D$Row").Value = $PlayerId
    $ws.Range("E$Row").Value = $PlayerName
    $ws.Range("F$Row").Value = $ProTeam
    $ws.Range("G$Row").Value = $LineupSlot
    # Force the nba_player_id to stay text (not auto-converted to a number)
    # while keeping the cell's original (default) style.
    $ws.Range("H$Row").Value = "'" + $NbaPlayerId
    $ws.Range("H$Row").Style = "Normal"
}

# Team KOBE: Kelly Oubre Jr. dropped, VJ Edgecombe added -> rows 12-14 shift up, row 15 is new
Set-RosterRow 12 4433247 "Jonathan Kuminga" "GSW" "PF" "1630228"
Set-RosterRow 13 3907497 "Dejounte Murray" "NOP" "SG" "1627749"
Set-RosterRow 14 3448 "Brook Lopez" "LAC" "C" "201572"
Set-RosterRow 15 5124612 "VJ Edgecombe" "PHL" "SG" "1642845"

# Team GNKI: Klay Thompson swapped out for Jalen Suggs
Set-RosterRow 29 4432165 "Jalen Suggs" "ORL" "PG" "1630591"

# Team PPH: Donte DiVincenzo dropped, Aaron Nesmith added -> row 124 shifts up, row 125 is new
Set-RosterRow 124 5037871 "Dylan Harper" "SAS" "PG" "1642844"
Set-RosterRow 125 4396909 "Aaron Nesmith" "IND" "SF" "1630174"

$wb.Save()
